$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the diff
$ws.Range("C2").Value = 12.2
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 11.1
$ws.Range("B4").Value = 0.7

# Update the active selection from C3 to C2
$ws.Range("C2").Select()
